$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '43.174.15'
$c.ClearFormats()
$ws.Range('E2').Value = '  +0.21%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.321.08'
$c.ClearFormats()
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  -0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '302.38'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.00%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '99.46'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.00%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.517'
$c.ClearFormats()
$ws.Range('E9').Value = '  +1.69%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.23'
$c.ClearFormats()
$ws.Range('E10').Value = '  +5.55%  '
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('E14').Value = '  +1.85%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '2.682.46'
$c.ClearFormats()
$ws.Range('E15').Value = '  +0.74%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.452.07'
$c.ClearFormats()
$ws.Range('E16').Value = '  +8.63%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.798'
$c.ClearFormats()
$ws.Range('E17').Value = '  -1.47%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '43.093.94'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('E21').Value = '  +0.11%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '68.27'
$c.ClearFormats()
$ws.Range('E22').Value = '  +0.59%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '240.56'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.41%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.15'
$c.ClearFormats()
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -0.03%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '25.49'
$c.ClearFormats()
$ws.Range('E27').Value = '  +3.60%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '168.93'
$c.ClearFormats()
$ws.Range('E28').Value = '  +1.58%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '34.40'
$c.ClearFormats()
$ws.Range('E29').Value = '  +1.30%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.20'
$c.ClearFormats()
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  -10.65%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.ClearFormats()
$ws.Range('E32').Value = '  +3.17%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E33').Value = '  -0.10%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.69'
$c.ClearFormats()
$ws.Range('E34').Value = '  +2.96%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '17.70'
$c.ClearFormats()
$ws.Range('E35').Value = '  +3.62%  '
$ws.Range('E36').Value = '  -0.98%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0696'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('E39').Value = '  +0.44%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('E41').Value = '  +0.08%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.004.73'
$c.ClearFormats()
$ws.Range('E42').Value = '  +0.33%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.0290'
$c.ClearFormats()
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('E44').Value = '  -4.77%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '10.12'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.33%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '17.68'
$c.ClearFormats()
$ws.Range('E46').Value = '  +0.01%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.89'
$c.ClearFormats()
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('E48').Value = '  +0.36%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '75.82'
$c.ClearFormats()
$ws.Range('E49').Value = '  +7.88%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.548.40'
$c.ClearFormats()
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('E51').Value = '  +1.82%  '
